$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the "Acme" reporting-mark row (col A) among the data rows and
# rename it to its correctly-cased form, "ACME".
$found = $ws.Range("A2:A126").Find("Acme")
if ($found -ne $null) {
    $found.Value2 = "ACME"
}

# Re-sort the data (A2:D126) alphabetically by Reporting Mark (col A),
# leaving the header row (row 1) untouched.
$dataRange = $ws.Range("A2:D126")
$dataRange.Sort($ws.Range("A2"))

# Match the saved selection state (active cell on the new ACME row).
[void]$ws.Range("A4").Select()
